# Generate Report for Handback
# Refresh the timestamp columns that record when handoff/handback xliff
# files were generated, as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for 6218c9b3-...md on the Overview sheet,
# and the matching "Correspond Handoff Datetime" for the same file on the
# de-de sheet, share the same text -- keep them in sync.
$wsOverview.Range("G2").Value = "2016-08-24 21:07:47"
$wsDeDe.Range("H2").Value = "2016-08-24 21:07:47"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the same row.
$wsZhCn.Range("H2").Value = "2016-08-24 21:07:42"
$wsZhCn.Range("K2").Value = "2016-08-24 21:07:59"

# de-de sheet: Correspond Handback DateTime.
$wsDeDe.Range("K2").Value = "2016-08-24 21:08:14"
